$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A new daily price record was inserted at row 210 (the dataset is ordered
# newest-first in blocks), which pushes the existing rows 210..343 down to
# 211..344. Insert a blank row at 210 to reproduce that shift.
$ws.Rows.Item(210).Insert()

# Populate the newly inserted row 210 with the new record's data.
$ws.Range("A210").Value = 10
$ws.Range("B210").Value = "Vega Modelo de Temuco"
$ws.Range("C210").Value = "La Araucanía"
$ws.Range("D210").Value = 44596
$ws.Range("E210").Value = 9
$ws.Range("F210").Value = 100114014
$ws.Range("G210").Value = "Betarraga"
$ws.Range("H210").Value = "Sin especificar"
$ws.Range("I210").Value = "Primera"
$ws.Range("J210").Value = 40
$ws.Range("K210").Value = 8000
$ws.Range("L210").Value = 8000
$ws.Range("M210").Value = 8000
$ws.Range("N210").Value = "$/docena de paquetes"
$ws.Range("O210").Value = "Provincia de Cautín"
$ws.Range("P210").Value = 667
$ws.Range("Q210").Value = 12
$ws.Range("R210").Value = "Hortaliza"
